$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("B2").Value = "https://www.headstruggle.de"
$ws.Range("D2").Value = "The Struggles in my Head"
$ws.Range("E2").Value = 24

# Delete row 3 entirely (the second URL entry is removed)
$ws.Rows("3:3").Delete()
